$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: label changes from "Sergipe" to "SE"
$ws.Range("A4").Value = "SE"

# Updated "Valor" (B column) figures
$ws.Range("B2").Value = 44.962263374485602
$ws.Range("B3").Value = 45.182469135802471
$ws.Range("B4").Value = 46.332222222222221
$ws.Range("B5").Value = 46.363333333333337
$ws.Range("B6").Value = 47.802222222222227
$ws.Range("B7").Value = 48.75333333333333
$ws.Range("B8").Value = 49.835555555555551
$ws.Range("B9").Value = 52.743333333333332
$ws.Range("B10").Value = 54.624444444444443

# Update selection to A2:C10 with active cell A2
$ws.Range("A2:C10").Select()
